$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A3").Value = "This is my first git branch work"
$ws.Range("A3").Select()
